$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'246.03"

# Row 3
$ws.Cells.Item(3, 4).Value = "'21.77"

# Row 4
$ws.Cells.Item(4, 4).Value = "'5.456"

# Row 5
$ws.Cells.Item(5, 4).Value = "'0.05660"

# Row 6
$ws.Cells.Item(6, 4).Value = "'3.377"

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.8025"

# Row 9
$ws.Cells.Item(9, 2).Value = 'WazirX'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(9, 4).Value = "'0.1433"
$ws.Cells.Item(9, 5).Value = '8WazirXWRX'

# Row 10
$ws.Cells.Item(10, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(10, 4).Value = "'0.07248"
$ws.Cells.Item(10, 5).Value = '9MandalaExchangeTokenMDX'

# Row 11
$ws.Cells.Item(11, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(11, 4).Value = "'0.03162"
$ws.Cells.Item(11, 5).Value = '10LiechtensteinCryptoassetsExchangeLCX'

# Row 12
$ws.Cells.Item(12, 2).Value = 'BitrueCoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12, 4).Value = "'0.02946"
$ws.Cells.Item(12, 5).Value = '11BitrueCoinBTR'

# Row 13
$ws.Cells.Item(13, 2).Value = 'BitMartToken'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(13, 4).Value = "'0.09284"
$ws.Cells.Item(13, 5).Value = '12BitMartTokenBMX'

# Row 14
$ws.Cells.Item(14, 2).Value = 'BitForexToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(14, 4).Value = "'0.001647"
$ws.Cells.Item(14, 5).Value = '13BitForexTokenBF'

# Row 15
$ws.Cells.Item(15, 2).Value = 'MCDex'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Cells.Item(15, 4).Value = "'3.212"
$ws.Cells.Item(15, 5).Value = '14MCDexMCB'

# Row 16
$ws.Cells.Item(16, 2).Value = 'CoinExToken'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Cells.Item(16, 4).Value = "'0.04724"
$ws.Cells.Item(16, 5).Value = '15CoinExTokenCET'

# Row 17
$ws.Cells.Item(17, 2).Value = 'One'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(17, 4).Value = "'0.0005898"
$ws.Cells.Item(17, 5).Value = '16OneONE'

# Row 18
$ws.Cells.Item(18, 4).Value = "'0.006383"

# Row 19
$ws.Cells.Item(19, 4).Value = "'0.005034"
$ws.Cells.Item(19, 5).Value = '18HotbitTokenHTBBestin24h'

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.0003204"

# Row 23
$ws.Cells.Item(23, 4).Value = "'3.810"

# Row 24
$ws.Cells.Item(24, 4).Value = "'6.430"

# Row 25
$ws.Cells.Item(25, 4).Value = "'2.126"

# Row 27
$ws.Cells.Item(27, 4).Value = "'0.1298"

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.04085"

# Row 41
$ws.Cells.Item(41, 2).Value = 'BKEXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Cells.Item(41, 4).Value = "'0.1041"
$ws.Cells.Item(41, 5).Value = '40BKEXTokenBKK'

# Row 42
$ws.Cells.Item(42, 4).Value = "'0.002973"

# Row 43
$ws.Cells.Item(43, 2).Value = 'KickToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Cells.Item(43, 4).Value = "'0.003242"
$ws.Cells.Item(43, 5).Value = '42KickTokenKICK'

# Row 44
$ws.Cells.Item(44, 4).Value = "'0.008068"

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.00005855"

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.00000000751"

# Row 47
$ws.Cells.Item(47, 4).Value = "'0.6834"

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.01051"

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.00002103"
